$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'29.934.83"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.05%  "

# Row 3
$ws.Range("D3").Value = "'1.895.70"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.13%  "

# Row 4
$ws.Range("D4").Value = "'1.000"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.03%  "

# Row 5
$ws.Range("D5").Value = "'0.7768"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.07%  "

# Row 6
$ws.Range("D6").Value = "'244.75"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.25%  "

# Row 7
$ws.Range("E7").Value = "  +0.01%  "

# Row 8
$ws.Range("D8").Value = "'0.3141"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.99%  "

# Row 9
$ws.Range("D9").Value = "'25.91"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.33%  "

# Row 10
$ws.Range("D10").Value = "'0.07272"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.09%  "

# Row 11
$ws.Range("D11").Value = "'0.09172"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +13.01%  "

# Row 12
$ws.Range("B12").Value = "Polygon"
$ws.Range("C12").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D12").Value = "'0.7755"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.63%  "

# Row 13
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "'1.920.59"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.01%  "

# Row 14
$ws.Range("D14").Value = "'5.453"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.36%  "

# Row 15
$ws.Range("D15").Value = "'94.85"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.25%  "

# Row 16
$ws.Range("D16").Value = "'6.234"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.96%  "

# Row 17
$ws.Range("D17").Value = "'29.947.57"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.09%  "

# Row 18
$ws.Range("D18").Value = "'14.00"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.02%  "

# Row 19
$ws.Range("D19").Value = "'247.05"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.63%  "

# Row 20
$ws.Range("D20").Value = "'0.000007905"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.59%  "

# Row 21
$ws.Range("D21").Value = "'2.165.63"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.17%  "

# Row 22
$ws.Range("D22").Value = "'8.151"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.78%  "

# Row 23
$ws.Range("D23").Value = "'1.000"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.10%  "

# Row 24
$ws.Range("D24").Value = "'1.000"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.01%  "

# Row 25
$ws.Range("D25").Value = "'0.1592"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -4.78%  "

# Row 26
$ws.Range("D26").Value = "'9.554"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.20%  "

# Row 27
$ws.Range("D27").Value = "'162.88"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.01%  "

# Row 28
$ws.Range("D28").Value = "'18.84"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.24%  "

# Row 29
$ws.Range("E29").Value = "  -1.26%  "

# Row 30
$ws.Range("D30").Value = "'1.424"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.39%  "

# Row 31
$ws.Range("E31").Value = "  +0.20%  "

# Row 32
$ws.Range("D32").Value = "'4.536"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.57%  "

# Row 33
$ws.Range("D33").Value = "'4.126"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.66%  "

# Row 34
$ws.Range("D34").Value = "'0.05524"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.31%  "

# Row 35
$ws.Range("D35").Value = "'1.251"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.77%  "

# Row 36
$ws.Range("D36").Value = "'0.7567"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.50%  "

# Row 37
$ws.Range("D37").Value = "'1.003"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.42%  "

# Row 38
$ws.Range("D38").Value = "'2.713"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.14%  "

# Row 39
$ws.Range("D39").Value = "'0.01979"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.96%  "

# Row 40
$ws.Range("E40").Value = "  +0.22%  "

# Row 41
$ws.Range("D41").Value = "'0.4513"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.62%  "

# Row 42
$ws.Range("D42").Value = "'74.39"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.37%  "

# Row 43
$ws.Range("D43").Value = "'6.106"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.44%  "

# Row 44
$ws.Range("D44").Value = "'1.094.80"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -6.29%  "

# Row 45
$ws.Range("D45").Value = "'0.8575"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.30%  "

# Row 46
$ws.Range("D46").Value = "'1.0000"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.01%  "

# Row 47
$ws.Range("D47").Value = "'1.896"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.03%  "

# Row 48
$ws.Range("D48").Value = "'102.82"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.85%  "

# Row 49
$ws.Range("D49").Value = "'7.636"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.67%  "

# Row 50
$ws.Range("D50").Value = "'9.889"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.29%  "

# Row 51
$ws.Range("D51").Value = "'3.009"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.23%  "
